$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "date" column (A2:A22) currently holds Excel date-serial values
# (31 Dec of each year, 2004-2024) rendered through a custom date number
# format. The naive forecaster component expects plain "YYYYQ4" text
# labels instead, so replace the values with text, and re-use the
# header cell's formatting (bold/bordered/centered, no custom date
# format) for the whole column.

# Copy the header cell's format (font, border, alignment, number
# format) onto the date column in one shot so the cells share the same
# style as the header instead of the old date-only style.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("A2:A22").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$startYear = 2004
$endYear = 2024

$row = 2
for ($year = $startYear; $year -le $endYear; $year++) {
    $ws.Cells.Item($row, 1).Value = "$($year)Q4"
    $row++
}
